$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

Replace-Text "2024-01-19 Friday" "2024-01-20 Saturday"
Replace-Text "90×38=3420" "71×24=1704"
Replace-Text "93×21=1953" "60×63=3780"
Replace-Text "83×24=1992" "82×54=4428"
Replace-Text "63×81=5103" "19×22=418"
Replace-Text "13×72=936" "71×62=4402"
Replace-Text "64×21=1344" "26×65=1690"
Replace-Text "89×42=3738" "86×80=6880"
Replace-Text "83×43=3569" "87×88=7656"
Replace-Text "70×97=6790" "44×13=572"
Replace-Text "65×44=2860" "77×60=4620"
Replace-Text "74×24=1776" "63×15=945"
Replace-Text "28×75=2100" "64×57=3648"
Replace-Text "32×84=2688" "68×11=748"
Replace-Text "15×98=1470" "24×91=2184"
Replace-Text "31×23=713" "22×81=1782"
Replace-Text "81×68=5508" "42×27=1134"
Replace-Text "49×67=3283" "36×86=3096"
Replace-Text "24×87=2088" "45×26=1170"
Replace-Text "39×81=3159" "64×33=2112"
Replace-Text "89×55=4895" "21×54=1134"
Replace-Text "89×20=1780" "17×42=714"
Replace-Text "32×60=1920" "62×68=4216"
Replace-Text "84×30=2520" "24×45=1080"
Replace-Text "91×21=1911" "86×45=3870"
Replace-Text "66×98=6468" "92×13=1196"
